$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title in A1
$ws.Range("A1").Value = "RESUMEN DE PEDIDOS DE COMPRA - CONSOLIDADO - VIVEVERDE 2026"

# Update row 9 (mascotas_manufacturado) values
$ws.Range("H9").Value = 160
$ws.Range("I9").Value = 60
$ws.Range("J9").Value = 2052
$ws.Range("K9").Value = 37.8
$ws.Range("L9").Value = 37
$ws.Range("Q9").Value = 66
